$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93-139 down to 94-140
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with new data
$ws.Cells.Item(93, 1).Value = 5
$ws.Cells.Item(93, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(93, 3).Value = "Maule"
$ws.Cells.Item(93, 4).Value = 44455
$ws.Cells.Item(93, 5).Value = 7
$ws.Cells.Item(93, 6).Value = 100112008
$ws.Cells.Item(93, 7).Value = "Coliflor"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 3000
$ws.Cells.Item(93, 11).Value = 600
$ws.Cells.Item(93, 12).Value = 600
$ws.Cells.Item(93, 13).Value = 600
$ws.Cells.Item(93, 14).Value = "`$/unidad"
$ws.Cells.Item(93, 15).Value = "Región del Maule"
$ws.Cells.Item(93, 16).Value = 600
$ws.Cells.Item(93, 17).Value = 1
$ws.Cells.Item(93, 18).Value = "Hortaliza"
